$d = $word.ActiveDocument

# 1. Remove the stray "_GoBack" bookmark (bookmarkStart/bookmarkEnd) left over
#    from the previous editing session.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Drop two bullet items from the list of methods - their wording is no
#    longer needed because the remaining bullets below them already cover
#    the same ground once they shift up:
#      - "Fazer a alocação de um bloco de disco ..."
#      - "Eliminar um ficheiro ou um diretório vazio pelo seu nome;"
$targets = @(
    "Fazer a alocação de um bloco de disco de forma a poder ser ocupado por um ficheiro (as pastas não ocupam espaço);",
    "Eliminar um ficheiro ou um diretório vazio pelo seu nome;"
)

foreach ($target in $targets) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Trim() -eq $target) {
            $p.Range.Delete()
            break
        }
    }
}
